$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row updates (E1/F1 renamed) ---
$ws.Range("E1").Value = "memories_add"
$ws.Range("F1").Value = "quantity_addmem"

# --- Row 1 height grows (wrap of new header text) ---
$ws.Rows.Item(1).RowHeight = 43.25

# --- New data added to row 2 (previously blank E2/F2) ---
$ws.Range("E2").Value = "nanya_pc2-4200U-444-12-A1"
$ws.Range("F2").Value = 6

# --- New column F width (narrow numeric column) ---
$ws.Columns.Item(6).ColumnWidth = 10.71

# --- New row 11 with a new processor / memory entry ---
$ws.Range("A11").Value = "usuario1@tmp.com"
$ws.Range("B11").Value = "7SF29N2"
$ws.Range("C11").Value = "00:00:00:00:00:10"
$ws.Range("E11").Value = "Fjeoeoe-eoeoeoeo"
$ws.Range("F11").Value = 4

# Mail hyperlink for the new user cell, matching the other rows' style
$ws.Hyperlinks.Add($ws.Range("A11"), "mailto:usuario1@tmp.com", [Type]::Missing, [Type]::Missing, "usuario1@tmp.com")

# Re-apply the same look-and-feel the other data rows use (avoid the auto
# "Hyperlink" theme style so row 11 matches rows 2..10 formatting)
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B10:C10").Copy()
$ws.Range("B11:C11").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E2").Copy()
$ws.Range("E11").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F2").Copy()
$ws.Range("F11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Keep the active selection pointing at E2 like the updated workbook ---
$ws.Range("E2").Select()
